$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename unit "GBX" -> "GBP" for Tesco and Sainsbury rows
$ws.Range("D2").Value = "GBP"
$ws.Range("D3").Value = "GBP"

# Update the active selection to D6, matching the saved view state
$ws.Range("D6").Select()
